$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 5352
$ws1.Range("F5").Value = 178
$ws1.Range("F6").Value = 223
$ws1.Range("F7").Value = 187
$ws1.Range("F8").Value = 8958
$ws1.Range("F10").Value = 649
$ws1.Range("F11").Value = 16
$ws1.Range("F12").Value = 2634
$ws1.Range("F13").Value = 2634
$ws1.Range("F14").Value = 6349
$ws1.Range("F15").Value = 2349
$ws1.Range("F19").Value = 2552
$ws1.Range("F21").Value = 22
$ws1.Range("F22").Value = 6632
$ws1.Range("F23").Value = 228
$ws1.Range("F25").Value = 160
$ws1.Range("F26").Value = 105
$ws1.Range("F28").Value = 7307
$ws1.Range("F32").Value = 46
$ws1.Range("F36").Value = 28
$ws1.Range("F39").Value = 59
$ws1.Range("F40").Value = 2557
$ws1.Range("F44").Value = 1140
$ws1.Range("F45").Value = 73
$ws1.Range("F46").Value = 566
$ws1.Range("F47").Value = 3579
$ws1.Range("F48").Value = 108
$ws1.Range("F49").Value = 1147
$ws1.Range("F50").Value = 45

# Sheet 2: 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 28
$ws2.Range("F5").Value = 219
$ws2.Range("F7").Value = 102
$ws2.Range("F8").Value = 20
$ws2.Range("F17").Value = 30

# Sheet 4: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5352
$ws4.Range("F4").Value = 5352
$ws4.Range("F5").Value = 178
$ws4.Range("F6").Value = 223
$ws4.Range("F7").Value = 8958
$ws4.Range("F9").Value = 649
$ws4.Range("F10").Value = 28
$ws4.Range("F11").Value = 2634
$ws4.Range("F14").Value = 219
$ws4.Range("F15").Value = 6349
$ws4.Range("F16").Value = 102
$ws4.Range("F20").Value = 2552
$ws4.Range("F22").Value = 20
$ws4.Range("F23").Value = 22
$ws4.Range("F24").Value = 6633
$ws4.Range("F25").Value = 228
$ws4.Range("F28").Value = 160
$ws4.Range("F29").Value = 105
$ws4.Range("F31").Value = 7307
$ws4.Range("F34").Value = 46
$ws4.Range("F40").Value = 59
$ws4.Range("F44").Value = 1140
$ws4.Range("F45").Value = 73
$ws4.Range("F46").Value = 3579
$ws4.Range("F47").Value = 108
$ws4.Range("F49").Value = 1147
$ws4.Range("F50").Value = 30
$ws4.Range("F51").Value = 45

